{"js": "// Locate the Conclusion / Recommendations / References heading paragraphs\n// by bookmark name, then insert the new body paragraphs described by the\n// diff: one FirstParagraph after \"Conclusion\", and a FirstParagraph +\n// BodyText pair after \"Recommendations\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst conclusionText =\n  \"Mr. Ramanujan went through a lot of hardships in his short lived life. \" +\n  \"Despite all of these things, he was able to pursue what he believed and \" +\n  \"earned a place amongst the greatest mathematicians in history. Being poor \" +\n  \"and having no degree, his persistence and ego has led him to find a job \" +\n  \"as a clerk. He then broke their tradition as Brahmins by cutting his hair \" +\n  \"and travelling overseas to publish his ideas. In Cambridge Trinity \" +\n  \"College, he worked with Mr. Hardy who taught him the rigor of proofs \" +\n  \"which led to the publication of his formulae and his becoming a Fellow \" +\n  \"of the Royal Society and a Fellow of Trinity College.\";\n\nconst recommendationsText1 =\n  \"The movie has made me appreciate pure mathematics. I have always evaded \" +\n  \"pure mathematics because of its abstract nature and the rigor of \" +\n  \"understanding and writing proofs. Nonetheless, I have started reading \" +\n  \"books on mathematical proofs and I am planning to delve into pure \" +\n  \"mathematics; an effect of the movie to me.\";\n\nconst recommendationsText2 =\n  \"It has also made me realize and appreciate the importance of HPS in \" +\n  \"teaching and learning mathematics. The history and philosophy \" +\n  \"mathematics that can be seen in the life story of Mr. Ramanujan and \" +\n  \"that of other mathematicians can spur appreciation and deeper \" +\n  \"understanding of mathematics. I would recommend my students to watch \" +\n  \"the film. The film showed how important proofs are to formulae, just \" +\n  \"like how critical thinking is important to science.\";\n\n// Find the Heading1 paragraphs that carry the section bookmarks by their\n// text, in document order (the body has exactly one of each).\nlet conclusionHeading = null;\nlet recommendationsHeading = null;\nfor (const p of paragraphs.items) {\n  p.load(\"text,style\");\n}\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  if (p.style === \"Heading 1\" && p.text.trim() === \"Conclusion\") {\n    conclusionHeading = p;\n  } else if (p.style === \"Heading 1\" && p.text.trim() === \"Recommendations\") {\n    recommendationsHeading = p;\n  }\n}\n\nif (!conclusionHeading || !recommendationsHeading) {\n  throw new Error(\"Could not locate Conclusion/Recommendations headings.\");\n}\n\n// Insert the Conclusion body paragraph right after the \"Conclusion\" heading.\nconst conclusionPara = conclusionHeading.insertParagraph(conclusionText, \"After\");\nconclusionPara.style = \"FirstParagraph\";\n\n// Insert the two Recommendations body paragraphs right after the\n// \"Recommendations\" heading, preserving their order.\nconst recPara1 = recommendationsHeading.insertParagraph(recommendationsText1, \"After\");\nrecPara1.style = \"FirstParagraph\";\n\nconst recPara2 = recPara1.insertParagraph(recommendationsText2, \"After\");\nrecPara2.style = \"BodyText\";\n\nawait context.sync();\n", "ps1": "# Insert the concluding paragraph after the \"Conclusion\" heading, and the\n# two recommendation paragraphs after the \"Recommendations\" heading, per\n# the commit \"finalizing CP_Ramanujan conclusion and recommendation\".\n\n$d = $word.ActiveDocument\n\n$conclusionText = \"Mr. Ramanujan went through a lot of hardships in his short lived life. \" +\n    \"Despite all of these things, he was able to pursue what he believed and \" +\n    \"earned a place amongst the greatest mathematicians in history. Being poor \" +\n    \"and having no degree, his persistence and ego has led him to find a job \" +\n    \"as a clerk. He then broke their tradition as Brahmins by cutting his hair \" +\n    \"and travelling overseas to publish his ideas. In Cambridge Trinity \" +\n    \"College, he worked with Mr. Hardy who taught him the rigor of proofs \" +\n    \"which led to the publication of his formulae and his becoming a Fellow \" +\n    \"of the Royal Society and a Fellow of Trinity College.\"\n\n$recommendationsText1 = \"The movie has made me appreciate pure mathematics. I have always evaded \" +\n    \"pure mathematics because of its abstract nature and the rigor of \" +\n    \"understanding and writing proofs. Nonetheless, I have started reading \" +\n    \"books on mathematical proofs and I am planning to delve into pure \" +\n    \"mathematics; an effect of the movie to me.\"\n\n$recommendationsText2 = \"It has also made me realize and appreciate the importance of HPS in \" +\n    \"teaching and learning mathematics. The history and philosophy \" +\n    \"mathematics that can be seen in the life story of Mr. Ramanujan and \" +\n    \"that of other mathematicians can spur appreciation and deeper \" +\n    \"understanding of mathematics. I would recommend my students to watch \" +\n    \"the film. The film showed how important proofs are to formulae, just \" +\n    \"like how critical thinking is important to science.\"\n\nfunction Get-HeadingParagraph($doc, $headingText) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.Trim() -eq $headingText) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# New body paragraph right after \"Conclusion\". Look the heading up fresh\n# (paragraph handles returned by the COM layer are position-anchored, not\n# logical, so anything fetched before a nearby insert can go stale).\n$conclusionHeading = Get-HeadingParagraph $d \"Conclusion\"\n$conclusionHeading.Range.InsertParagraphAfter()\n$conclusionPara = $conclusionHeading.Next()\n$conclusionPara.Range.Text = $conclusionText\n$conclusionPara.Style = \"FirstParagraph\"\n\n# Two new body paragraphs right after \"Recommendations\", in order. Re-fetch\n# the heading now that the document has grown from the edit above.\n$recommendationsHeading = Get-HeadingParagraph $d \"Recommendations\"\n$recommendationsHeading.Range.InsertParagraphAfter()\n$recPara1 = $recommendationsHeading.Next()\n$recPara1.Range.Text = $recommendationsText1\n$recPara1.Style = \"FirstParagraph\"\n\n$recPara1.Range.InsertParagraphAfter()\n$recPara2 = $recPara1.Next()\n$recPara2.Range.Text = $recommendationsText2\n$recPara2.Style = \"BodyText\"\n"}
